# ContosoLearn Market Research - translate body text from Chinese to English
# and normalize paragraph/run formatting (drop the zh-CN rPr bag and the
# explicit w:bidi="0" paragraph property) to match the authored edit.

$d = $word.ActiveDocument

function Set-ParaBodyXml($Index, $InnerXml) {
    $p = $d.Paragraphs($Index)
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $InnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $p.Range.InsertXML($pkg)
}

$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr>'

# 1. Title
$x1 = '<w:p><w:r><w:t>ContosoLearn Market Research</w:t></w:r></w:p>'
Set-ParaBodyXml 1 $x1

# 2. AdatumLearn (two source runs merge into one)
$x2 = '<w:p>' + $listPPr + '<w:r><w:t>AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology.</w:t></w:r></w:p>'
Set-ParaBodyXml 2 $x2

# 3. AdventureLearn
$x3 = '<w:p>' + $listPPr + '<w:r><w:t>AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations.</w:t></w:r></w:p>'
Set-ParaBodyXml 3 $x3

# 4. AlpineTraining
$x4 = '<w:p>' + $listPPr + '<w:r><w:t>AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning.</w:t></w:r></w:p>'
Set-ParaBodyXml 4 $x4

# 5. Bellows OnDemand
$x5 = '<w:p>' + $listPPr + '<w:r><w:t>Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration.</w:t></w:r></w:p>'
Set-ParaBodyXml 5 $x5

# 6. FabrikamLearning
$x6 = '<w:p>' + $listPPr + '<w:r><w:t>FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs.</w:t></w:r></w:p>'
Set-ParaBodyXml 6 $x6

# 7. FirstUp Cards
$x7 = '<w:p>' + $listPPr + '<w:r><w:t>FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario.</w:t></w:r></w:p>'
Set-ParaBodyXml 7 $x7

# 8. Munson'sLearn
$x8 = '<w:p>' + $listPPr + "<w:r><w:t>Munson'sLearn: Munson'sLearn is designed to enable businesses to train their employees, partners, and customers.</w:t></w:r></w:p>"
Set-ParaBodyXml 8 $x8

# 9. LibertyLearn
$x9 = '<w:p>' + $listPPr + '<w:r><w:t>LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project.</w:t></w:r></w:p>'
Set-ParaBodyXml 9 $x9

# 10. WoodgroveLMS - split across three runs with a gramStart/gramEnd proofErr pair
$x10 = '<w:p>' + $listPPr + '<w:r><w:t xml:space="preserve">WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a best</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>-in-class training experience.</w:t></w:r></w:p>'
Set-ParaBodyXml 10 $x10

# 11. NorthwindWorlds
$x11 = '<w:p>' + $listPPr + '<w:r><w:t>NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises.</w:t></w:r></w:p>'
Set-ParaBodyXml 11 $x11

# 12. ProsewareLearn
$x12 = '<w:p>' + $listPPr + '<w:r><w:t>ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website.</w:t></w:r></w:p>'
Set-ParaBodyXml 12 $x12

# 13. RelecloudLearn
$x13 = '<w:p>' + $listPPr + '<w:r><w:t>RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects.</w:t></w:r></w:p>'
Set-ParaBodyXml 13 $x13

# 14. TreyAcademy
$x14 = '<w:p>' + $listPPr + '<w:r><w:t>TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010.</w:t></w:r></w:p>'
Set-ParaBodyXml 14 $x14

# 15. Closing paragraph - three source runs merge into a single run
$x15 = '<w:p><w:r><w:t xml:space="preserve">These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. </w:t></w:r></w:p>'
Set-ParaBodyXml 15 $x15
